$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed odds in existing rows 3, 4, 7 ---
# Row 3 changes
$ws.Cells.Item(3, 7).Value = 2.7
$ws.Cells.Item(3, 9).Value = 2.55
$ws.Cells.Item(3, 26).Value = 26
$ws.Cells.Item(3, 36).Value = 10
$ws.Cells.Item(3, 37).Value = 26
$ws.Cells.Item(3, 43).Value = 41
$ws.Cells.Item(3, 50).Value = 4.75
$ws.Cells.Item(3, 55).Value = 151

# Row 4 changes
$ws.Cells.Item(4, 16).Value = 5.55
$ws.Cells.Item(4, 17).Value = 1.33
$ws.Cells.Item(4, 18).Value = 2.73

# Row 7 changes
$ws.Cells.Item(7, 7).Value = 3.4
$ws.Cells.Item(7, 8).Value = 3.3
$ws.Cells.Item(7, 9).Value = 1.95
$ws.Cells.Item(7, 10).Value = 3.85
$ws.Cells.Item(7, 11).Value = 2.12
$ws.Cells.Item(7, 13).Value = 10.4
$ws.Cells.Item(7, 14).Value = 1.01
$ws.Cells.Item(7, 15).Value = 1.18
$ws.Cells.Item(7, 16).Value = 3.66
$ws.Cells.Item(7, 17).Value = 1.7
$ws.Cells.Item(7, 18).Value = 1.93
$ws.Cells.Item(7, 19).Value = 1.35
$ws.Cells.Item(7, 20).Value = 3.04
$ws.Cells.Item(7, 21).Value = 1.61
$ws.Cells.Item(7, 22).Value = 2.25
$ws.Cells.Item(7, 23).Value = 9.5
$ws.Cells.Item(7, 24).Value = 16
$ws.Cells.Item(7, 27).Value = 23
$ws.Cells.Item(7, 28).Value = 25
$ws.Cells.Item(7, 29).Value = 11
$ws.Cells.Item(7, 30).Value = 5.8
$ws.Cells.Item(7, 31).Value = 10.5
$ws.Cells.Item(7, 32).Value = 37
$ws.Cells.Item(7, 33).Value = 200
$ws.Cells.Item(7, 34).Value = 7.3
$ws.Cells.Item(7, 35).Value = 9
$ws.Cells.Item(7, 37).Value = 15
$ws.Cells.Item(7, 39).Value = 17.5
$ws.Cells.Item(7, 41).Value = 19
$ws.Cells.Item(7, 42).Value = 23
$ws.Cells.Item(7, 43).Value = 90
$ws.Cells.Item(7, 44).Value = 120
$ws.Cells.Item(7, 46).Value = 2.87
$ws.Cells.Item(7, 47).Value = 6.6
$ws.Cells.Item(7, 48).Value = 50
$ws.Cells.Item(7, 50).Value = 4

# --- Insert two new rows at position 9 (shifts old rows 9,10 -> 11,12) ---
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# Fill new row 9 (data)
$ws.Cells.Item(9, 1).Value = "xdMjMez3"
$ws.Cells.Item(9, 2).Value = "26/10/2024"
$ws.Cells.Item(9, 3).Value = "07:00"
$ws.Cells.Item(9, 4).Value = "PORTUGAL - LIGA PORTUGAL 2"
$ws.Cells.Item(9, 5).Value = "Penafiel"
$ws.Cells.Item(9, 6).Value = "Academico Viseu"
$ws.Cells.Item(9, 7).Value = 2.25
$ws.Cells.Item(9, 8).Value = 3.1
$ws.Cells.Item(9, 9).Value = 3.1
$ws.Cells.Item(9, 10).Value = 3
$ws.Cells.Item(9, 11).Value = 2.1
$ws.Cells.Item(9, 12).Value = 3.75
$ws.Cells.Item(9, 13).Value = 1.06
$ws.Cells.Item(9, 14).Value = 10
$ws.Cells.Item(9, 15).Value = 1.33
$ws.Cells.Item(9, 16).Value = 3.25
$ws.Cells.Item(9, 17).Value = 2.05
$ws.Cells.Item(9, 18).Value = 1.75
$ws.Cells.Item(9, 19).Value = 1.44
$ws.Cells.Item(9, 20).Value = 2.63
$ws.Cells.Item(9, 21).Value = 1.8
$ws.Cells.Item(9, 22).Value = 1.91
$ws.Cells.Item(9, 23).Value = 8
$ws.Cells.Item(9, 24).Value = 11
$ws.Cells.Item(9, 25).Value = 9.5
$ws.Cells.Item(9, 26).Value = 21
$ws.Cells.Item(9, 27).Value = 19
$ws.Cells.Item(9, 28).Value = 29
$ws.Cells.Item(9, 29).Value = 9
$ws.Cells.Item(9, 30).Value = 6
$ws.Cells.Item(9, 31).Value = 13
$ws.Cells.Item(9, 32).Value = 51
$ws.Cells.Item(9, 33).Value = 251
$ws.Cells.Item(9, 34).Value = 9.5
$ws.Cells.Item(9, 35).Value = 15
$ws.Cells.Item(9, 36).Value = 12
$ws.Cells.Item(9, 37).Value = 34
$ws.Cells.Item(9, 38).Value = 26
$ws.Cells.Item(9, 39).Value = 34
$ws.Cells.Item(9, 40).Value = 4.33
$ws.Cells.Item(9, 41).Value = 13
$ws.Cells.Item(9, 42).Value = 23
$ws.Cells.Item(9, 43).Value = 41
$ws.Cells.Item(9, 44).Value = 67
$ws.Cells.Item(9, 45).Value = 151
$ws.Cells.Item(9, 46).Value = 2.63
$ws.Cells.Item(9, 47).Value = 8
$ws.Cells.Item(9, 48).Value = 51
$ws.Cells.Item(9, 49).Value = 81
$ws.Cells.Item(9, 50).Value = 5
$ws.Cells.Item(9, 51).Value = 17
$ws.Cells.Item(9, 52).Value = 26
$ws.Cells.Item(9, 53).Value = 51
$ws.Cells.Item(9, 54).Value = 81
$ws.Cells.Item(9, 55).Value = 201
$ws.Cells.Item(9, 56).Value = 81

# Fill new row 10 (data)
$ws.Cells.Item(10, 1).Value = "p4hZVNR1"
$ws.Cells.Item(10, 2).Value = "26/10/2024"
$ws.Cells.Item(10, 3).Value = "07:30"
$ws.Cells.Item(10, 4).Value = "ROMANIA - LIGA 1"
$ws.Cells.Item(10, 5).Value = "UTA Arad"
$ws.Cells.Item(10, 6).Value = "Gloria Buzau"
$ws.Cells.Item(10, 7).Value = 1.83
$ws.Cells.Item(10, 8).Value = 3.3
$ws.Cells.Item(10, 9).Value = 4.33
$ws.Cells.Item(10, 10).Value = 2.5
$ws.Cells.Item(10, 11).Value = 2.1
$ws.Cells.Item(10, 12).Value = 4.5
$ws.Cells.Item(10, 13).Value = 1.07
$ws.Cells.Item(10, 14).Value = 9
$ws.Cells.Item(10, 15).Value = 1.3
$ws.Cells.Item(10, 16).Value = 3.4
$ws.Cells.Item(10, 17).Value = 2.05
$ws.Cells.Item(10, 18).Value = 1.8
$ws.Cells.Item(10, 19).Value = 1.4
$ws.Cells.Item(10, 20).Value = 2.75
$ws.Cells.Item(10, 21).Value = 1.83
$ws.Cells.Item(10, 22).Value = 1.83
$ws.Cells.Item(10, 23).Value = 7
$ws.Cells.Item(10, 24).Value = 8.5
$ws.Cells.Item(10, 25).Value = 8.5
$ws.Cells.Item(10, 26).Value = 15
$ws.Cells.Item(10, 27).Value = 15
$ws.Cells.Item(10, 28).Value = 29
$ws.Cells.Item(10, 29).Value = 9
$ws.Cells.Item(10, 30).Value = 6.5
$ws.Cells.Item(10, 31).Value = 15
$ws.Cells.Item(10, 32).Value = 51
$ws.Cells.Item(10, 33).Value = 301
$ws.Cells.Item(10, 34).Value = 12
$ws.Cells.Item(10, 35).Value = 21
$ws.Cells.Item(10, 36).Value = 15
$ws.Cells.Item(10, 37).Value = 41
$ws.Cells.Item(10, 38).Value = 34
$ws.Cells.Item(10, 39).Value = 41
$ws.Cells.Item(10, 40).Value = 3.75
$ws.Cells.Item(10, 41).Value = 10
$ws.Cells.Item(10, 42).Value = 21
$ws.Cells.Item(10, 43).Value = 34
$ws.Cells.Item(10, 44).Value = 51
$ws.Cells.Item(10, 45).Value = 151
$ws.Cells.Item(10, 46).Value = 2.75
$ws.Cells.Item(10, 47).Value = 8.5
$ws.Cells.Item(10, 48).Value = 51
$ws.Cells.Item(10, 49).Value = 51
$ws.Cells.Item(10, 50).Value = 6
$ws.Cells.Item(10, 51).Value = 23
$ws.Cells.Item(10, 52).Value = 29
$ws.Cells.Item(10, 53).Value = 81
$ws.Cells.Item(10, 54).Value = 101
$ws.Cells.Item(10, 55).Value = 251
$ws.Cells.Item(10, 56).Value = 51

# --- New rows 13, 14 appended after row 12 (no insert needed; writing past the end extends the sheet) ---
# Fill new row 13 (data)
$ws.Cells.Item(13, 1).Value = "zocXixyg"
$ws.Cells.Item(13, 2).Value = "26/10/2024"
$ws.Cells.Item(13, 3).Value = "07:30"
$ws.Cells.Item(13, 4).Value = "TURKEY - SUPER LIG"
$ws.Cells.Item(13, 5).Value = "Hatayspor"
$ws.Cells.Item(13, 6).Value = "Kayserispor"
$ws.Cells.Item(13, 7).Value = 2.2
$ws.Cells.Item(13, 8).Value = 3.5
$ws.Cells.Item(13, 9).Value = 3.1
$ws.Cells.Item(13, 10).Value = 2.88
$ws.Cells.Item(13, 11).Value = 2.25
$ws.Cells.Item(13, 12).Value = 3.5
$ws.Cells.Item(13, 13).Value = 1.04
$ws.Cells.Item(13, 14).Value = 13
$ws.Cells.Item(13, 15).Value = 1.22
$ws.Cells.Item(13, 16).Value = 4
$ws.Cells.Item(13, 17).Value = 1.73
$ws.Cells.Item(13, 18).Value = 2.08
$ws.Cells.Item(13, 19).Value = 1.33
$ws.Cells.Item(13, 20).Value = 3.25
$ws.Cells.Item(13, 21).Value = 1.62
$ws.Cells.Item(13, 22).Value = 2.2
$ws.Cells.Item(13, 23).Value = 9.5
$ws.Cells.Item(13, 24).Value = 12
$ws.Cells.Item(13, 25).Value = 9
$ws.Cells.Item(13, 26).Value = 21
$ws.Cells.Item(13, 27).Value = 17
$ws.Cells.Item(13, 28).Value = 23
$ws.Cells.Item(13, 29).Value = 13
$ws.Cells.Item(13, 30).Value = 7
$ws.Cells.Item(13, 31).Value = 13
$ws.Cells.Item(13, 32).Value = 41
$ws.Cells.Item(13, 33).Value = 151
$ws.Cells.Item(13, 34).Value = 12
$ws.Cells.Item(13, 35).Value = 17
$ws.Cells.Item(13, 36).Value = 11
$ws.Cells.Item(13, 37).Value = 34
$ws.Cells.Item(13, 38).Value = 23
$ws.Cells.Item(13, 39).Value = 29
$ws.Cells.Item(13, 40).Value = 4.5
$ws.Cells.Item(13, 41).Value = 12
$ws.Cells.Item(13, 42).Value = 19
$ws.Cells.Item(13, 43).Value = 41
$ws.Cells.Item(13, 44).Value = 51
$ws.Cells.Item(13, 45).Value = 126
$ws.Cells.Item(13, 46).Value = 3.25
$ws.Cells.Item(13, 47).Value = 7.5
$ws.Cells.Item(13, 48).Value = 41
$ws.Cells.Item(13, 49).Value = 251
$ws.Cells.Item(13, 50).Value = 5
$ws.Cells.Item(13, 51).Value = 17
$ws.Cells.Item(13, 52).Value = 23
$ws.Cells.Item(13, 53).Value = 51
$ws.Cells.Item(13, 54).Value = 67
$ws.Cells.Item(13, 55).Value = 151
$ws.Cells.Item(13, 56).Value = 301

# Fill new row 14 (data)
$ws.Cells.Item(14, 1).Value = "8p07HlDr"
$ws.Cells.Item(14, 2).Value = "26/10/2024"
$ws.Cells.Item(14, 3).Value = "07:30"
$ws.Cells.Item(14, 4).Value = "TURKEY - 1. LIG"
$ws.Cells.Item(14, 5).Value = "Corum"
$ws.Cells.Item(14, 6).Value = "Genclerbirligi"
$ws.Cells.Item(14, 7).Value = 2.25
$ws.Cells.Item(14, 8).Value = 3.2
$ws.Cells.Item(14, 9).Value = 3
$ws.Cells.Item(14, 10).Value = 3.1
$ws.Cells.Item(14, 11).Value = 2.05
$ws.Cells.Item(14, 12).Value = 3.75
$ws.Cells.Item(14, 13).Value = 1.07
$ws.Cells.Item(14, 14).Value = 8.5
$ws.Cells.Item(14, 15).Value = 1.36
$ws.Cells.Item(14, 16).Value = 3
$ws.Cells.Item(14, 17).Value = 2.2
$ws.Cells.Item(14, 18).Value = 1.65
$ws.Cells.Item(14, 19).Value = 1.5
$ws.Cells.Item(14, 20).Value = 2.5
$ws.Cells.Item(14, 21).Value = 1.91
$ws.Cells.Item(14, 22).Value = 1.8
$ws.Cells.Item(14, 23).Value = 7
$ws.Cells.Item(14, 24).Value = 10
$ws.Cells.Item(14, 25).Value = 9.5
$ws.Cells.Item(14, 26).Value = 21
$ws.Cells.Item(14, 27).Value = 21
$ws.Cells.Item(14, 28).Value = 34
$ws.Cells.Item(14, 29).Value = 8.5
$ws.Cells.Item(14, 30).Value = 6.5
$ws.Cells.Item(14, 31).Value = 17
$ws.Cells.Item(14, 32).Value = 51
$ws.Cells.Item(14, 33).Value = 351
$ws.Cells.Item(14, 34).Value = 8.5
$ws.Cells.Item(14, 35).Value = 15
$ws.Cells.Item(14, 36).Value = 12
$ws.Cells.Item(14, 37).Value = 34
$ws.Cells.Item(14, 38).Value = 29
$ws.Cells.Item(14, 39).Value = 41
$ws.Cells.Item(14, 40).Value = 4.33
$ws.Cells.Item(14, 41).Value = 13
$ws.Cells.Item(14, 42).Value = 26
$ws.Cells.Item(14, 43).Value = 41
$ws.Cells.Item(14, 44).Value = 67
$ws.Cells.Item(14, 45).Value = 201
$ws.Cells.Item(14, 46).Value = 2.5
$ws.Cells.Item(14, 47).Value = 8.5
$ws.Cells.Item(14, 48).Value = 51
$ws.Cells.Item(14, 49).Value = 126
$ws.Cells.Item(14, 50).Value = 5
$ws.Cells.Item(14, 51).Value = 19
$ws.Cells.Item(14, 52).Value = 29
$ws.Cells.Item(14, 53).Value = 51
$ws.Cells.Item(14, 54).Value = 81
$ws.Cells.Item(14, 55).Value = 251
$ws.Cells.Item(14, 56).Value = 126
